# Update "人气值/浏览量" (column F) figures on the worksheets, reflecting
# freshly generated output data (gh-pages regeneration).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F9").Value  = 284
$ws1.Range("F10").Value = 407
$ws1.Range("F11").Value = 361
$ws1.Range("F12").Value = 1806
$ws1.Range("F13").Value = 791
$ws1.Range("F18").Value = 1291
$ws1.Range("F22").Value = 373
$ws1.Range("F25").Value = 119
$ws1.Range("F26").Value = 6780
$ws1.Range("F27").Value = 6354
$ws1.Range("F37").Value = 1328

# Sheet "本地生活" (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F5").Value = 87

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value  = 87
$ws4.Range("F12").Value = 284
$ws4.Range("F14").Value = 407
$ws4.Range("F15").Value = 361
$ws4.Range("F16").Value = 1806
$ws4.Range("F17").Value = 791
$ws4.Range("F22").Value = 1291
$ws4.Range("F24").Value = 373
$ws4.Range("F26").Value = 119
$ws4.Range("F29").Value = 6780
$ws4.Range("F30").Value = 6357
$ws4.Range("F33").Value = 1328
